# "add update quantity feature"
# Update the single data row (row 2) on the active sheet with the new
# shopping-trip details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text fields - safe to assign directly, Excel keeps these as text.
$ws.Range("A2").Value = "shopping with Diego Rivera"
$ws.Range("E2").Value = "Small Copper Bench"
$ws.Range("F2").Value = "Mung Beans"

# "1975-07-11" and "4.0" both *look* like a date / number, so a plain
# Range.Value assignment would make Excel auto-convert them (date serial /
# numeric value) instead of keeping them as the literal text that is stored
# in the workbook. To force them to stay text (matching the original
# cell type) without disturbing the existing cell style, stage the text in
# a scratch cell as a text formula, then copy/paste-special just the
# resulting value onto the target cell - this preserves the destination's
# existing formatting exactly and leaves no extra styles behind.

$scratch = $ws.Range("Z100")

$scratch.Formula = '="1975-07-11"'
$scratch.Copy()
$ws.Range("B2").PasteSpecial(-4163) # xlPasteValues
$scratch.Clear()

$scratch.Formula = '="4.0"'
$scratch.Copy()
$ws.Range("G2").PasteSpecial(-4163) # xlPasteValues
$scratch.Clear()
